$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")
$data.Range("A:A").Delete()

$source = $wb.Worksheets.Item("source")
foreach ($r in 2, 7, 10) {
    $a = $source.Range("A" + $r)
    $a.Borders.Item(7).LineStyle = 1
    $a.Borders.Item(7).Weight = 2
    $a.Borders.Item(8).LineStyle = 1
    $a.Borders.Item(8).Weight = 2
    $a.Borders.Item(9).LineStyle = 1
    $a.Borders.Item(9).Weight = 2
    $a.HorizontalAlignment = -4108
    $a.VerticalAlignment = -4108
    $a.WrapText = $true

    $b = $source.Range("B" + $r)
    $b.Borders.Item(10).LineStyle = 1
    $b.Borders.Item(10).Weight = 2
    $b.Borders.Item(8).LineStyle = 1
    $b.Borders.Item(8).Weight = 2
    $b.Borders.Item(9).LineStyle = 1
    $b.Borders.Item(9).Weight = 2
    $b.HorizontalAlignment = -4108
    $b.VerticalAlignment = -4108
    $b.WrapText = $true
}
